$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Simple text/string updates (not numeric-looking, safe to assign directly)
$ws.Range("D2").Value = "29.900.93"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.631.28"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("D12").Value = "1.866.06"
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("D13").Value = "1.633.50"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("E15").Value = "  +8.50%  "
$ws.Range("D16").Value = "29.916.21"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "0.0₃0701"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("E24").Value = "  +2.33%  "
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("E31").Value = "  +3.24%  "
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").Value = "1.424.54"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("E35").Value = "  +3.53%  "
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("E40").Value = "  +9.92%  "
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("E47").Value = "  -5.77%  "
$ws.Range("D48").Value = "1.773.66"
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("E49").Value = "  -2.23%  "
$ws.Range("E50").Value = "  +3.77%  "
$ws.Range("D51").Value = "0.0₆0113"
$ws.Range("E51").Value = "  +18.38%  "

# Numeric-looking text values in column D: force text type to preserve formatting
# (e.g. leading/trailing zeros) without leaving the cell as a real number.
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value2 = "214.56"
$r.Style = "Normal"
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value2 = "28.69"
$r.Style = "Normal"
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value2 = "0.562"
$r.Style = "Normal"
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value2 = "9.33"
$r.Style = "Normal"
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value2 = "64.10"
$r.Style = "Normal"
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value2 = "241.02"
$r.Style = "Normal"
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value2 = "9.80"
$r.Style = "Normal"
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value2 = "158.19"
$r.Style = "Normal"
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value2 = "0.110"
$r.Style = "Normal"
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value2 = "1.03"
$r.Style = "Normal"
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value2 = "75.53"
$r.Style = "Normal"
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value2 = "0.556"
$r.Style = "Normal"
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value2 = "0.827"
$r.Style = "Normal"
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value2 = "50.74"
$r.Style = "Normal"
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value2 = "90.44"
$r.Style = "Normal"
